$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $c = $ws.Range($cellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $newValue
    $c.Style = $origStyle
}

Set-TextValue "D2" "292.42"
Set-TextValue "E2" "-6.81%"
Set-TextValue "D3" "40.44"
Set-TextValue "E3" "0.20%"
Set-TextValue "D4" "5.027"
Set-TextValue "E4" "-2.70%"
Set-TextValue "D5" "0.07326"
Set-TextValue "E5" "-3.41%"
Set-TextValue "B6" "FTXToken"
Set-TextValue "C6" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D6" "1.533"
Set-TextValue "E6" "-7.90%"
Set-TextValue "B7" "MXToken"
Set-TextValue "C7" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D7" "0.9272"
Set-TextValue "E7" "-0.03%"
Set-TextValue "B8" "BTSEToken"
Set-TextValue "C8" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D8" "2.369"
Set-TextValue "E8" "-2.27%"
Set-TextValue "B9" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C9" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D9" "0.1185"
Set-TextValue "E9" "-1.14%"
Set-TextValue "B10" "WazirX"
Set-TextValue "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1744"
Set-TextValue "E10" "-4.17%"
Set-TextValue "B11" "BitrueCoin"
Set-TextValue "C11" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D11" "0.04328"
Set-TextValue "E11" "4.08%"
Set-TextValue "B12" "MandalaExchangeToken"
Set-TextValue "C12" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D12" "0.08645"
Set-TextValue "E12" "-4.37%"
Set-TextValue "B13" "BitMartToken"
Set-TextValue "C13" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D13" "0.1054"
Set-TextValue "E13" "0.16%"
Set-TextValue "B14" "BitForexToken"
Set-TextValue "C14" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D14" "0.001274"
Set-TextValue "E14" "-0.47%"
Set-TextValue "B15" "TigerCash"
Set-TextValue "C15" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D15" "0.005971"
Set-TextValue "E15" "2.92%"
Set-TextValue "B16" "LEO"
Set-TextValue "C16" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D16" "3.340"
Set-TextValue "E16" "0.26%"
Set-TextValue "B17" "GateToken"
Set-TextValue "C17" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D17" "4.300"
Set-TextValue "E17" "-0.63%"
Set-TextValue "E18" "-2.01%"
Set-TextValue "E19" "5.55%"
Set-TextValue "E20" "2.88%"
Set-TextValue "E22" "-2.14%"
Set-TextValue "D23" "0.001261"
Set-TextValue "E23" "-0.82%"
Set-TextValue "D24" "0.003781"
Set-TextValue "E25" "0.81%"
Set-TextValue "D38" "0.02275"
Set-TextValue "E38" "-5.65%"
Set-TextValue "D39" "0.04968"
Set-TextValue "E39" "-3.45%"
Set-TextValue "E40" "70.58%"
Set-TextValue "D41" "0.007695"
Set-TextValue "E41" "-0.30%"
Set-TextValue "D42" "0.1286"
Set-TextValue "E42" "-1.10%"
Set-TextValue "D43" "0.007361"
Set-TextValue "E43" "-3.10%"
Set-TextValue "D44" "0.008301"
Set-TextValue "E44" "-3.00%"
Set-TextValue "D45" "0.2917"
Set-TextValue "E45" "-14.02%"
Set-TextValue "D46" "0.00006307"
Set-TextValue "E46" "-4.30%"
Set-TextValue "D48" "0.03728"
Set-TextValue "E48" "-86.12%"
